# Boston roster: re-sort two adjacent player rows (swap all columns
# except the "No." index column A) so that:
#   - Row 6 (Luke Kornet) <-> Row 7 (Malcolm Brogdon)
#   - Row 14 (JD Davison (TW)) <-> Row 15 (Mike Muscala)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

# Column I ("Exp") stores values such as "5", "6", "9", "13" as TEXT
# (shared strings, because the column also contains "R" for rookies).
# A plain "$cell.Value2 = '6'" assignment lets the host auto-coerce a
# numeric-looking string into a real number, which would change the
# cell's type (t="s" -> numeric) relative to the original file. Forcing
# a text NumberFormat for the duration of the write keeps it textual;
# restoring the "Normal" style afterwards avoids leaving a stray
# number-format behind on the cell.
$textCols = @("I")

function Set-CellValue {
    param($cell, $value, $isTextCol)
    if ($isTextCol) {
        $cell.NumberFormat = "@"
        $cell.Value2 = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value2 = $value
    }
}

function Swap-Rows {
    param($ws, $row1, $row2, $cols, $textCols)
    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$row1")
        $cell2 = $ws.Range("$col$row2")
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $isTextCol = $textCols -contains $col
        Set-CellValue $cell1 $v2 $isTextCol
        Set-CellValue $cell2 $v1 $isTextCol
    }
}

Swap-Rows $ws 6 7 $cols $textCols
Swap-Rows $ws 14 15 $cols $textCols
